# The commit shuffles the *text* of several runs around the document
# (paragraph structure, styles and run formatting are untouched). Each
# of the 11 runs below gets a new value that, before the edit, lived in
# a *different* run elsewhere in the document - i.e. this is a cyclic
# rotation of text among fixed structural slots.
#
# Doing the replacements directly (old -> new) would cascade: a rule's
# "new" text can equal a later rule's "old" text, so that later rule
# would incorrectly re-match text that a previous rule just wrote.
# To avoid that we stage everything through unique placeholder tokens
# first, then fill in the final text in a second pass.

$d = $word.ActiveDocument

# Keep straight apostrophes/quotes as-is (don't let AutoFormat turn the
# apostrophe in "students'" into a curly quote on replace).
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function Replace-Text($old, $new) {
    $d.Content.Find.Execute(
        $old,   # FindText
        $true,  # MatchCase
        $false, # MatchWholeWord
        $false, # MatchWildcards
        $false, # MatchSoundsLike
        $false, # MatchAllWordForms
        $true,  # Forward
        1,      # Wrap (wdFindContinue)
        $false, # Format
        $new,   # ReplaceWith
        2       # Replace (wdReplaceAll)
    ) | Out-Null
}

# Phase 1: old text -> unique placeholder
Replace-Text "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte." "@@SLOT0@@"
Replace-Text "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics." "@@SLOT1@@"
Replace-Text "5840730 - Antonio Jefferson da Silva Machado" "@@SLOT2@@"
Replace-Text "519033 - Carlos Yujiro Shigue" "@@SLOT3@@"
Replace-Text "A definir, de acordo com o tópico programado." "@@SLOT4@@"
Replace-Text "To be defined, according to the programmed topic." "@@SLOT5@@"
Replace-Text "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação." "@@SLOT6@@"
Replace-Text "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa." "@@SLOT7@@"
Replace-Text "A média do semestre será computada com base na relação: M=(P1+2P2)/3" "@@SLOT8@@"
Replace-Text "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será computada com base na relação: MF=(M+RC)/2" "@@SLOT9@@"
Replace-Text "Apostila ou texto fornecido pelo docente responsável. Artigos extraídos de revistas especializadas nas áreas de Ciências e Tecnologia." "@@SLOT10@@"

# Phase 2: placeholder -> final text
Replace-Text "@@SLOT0@@" "A definir, de acordo com o tópico programado."
Replace-Text "@@SLOT1@@" "To be defined, according to the programmed topic."
Replace-Text "@@SLOT2@@" "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
Replace-Text "@@SLOT3@@" "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
Replace-Text "@@SLOT4@@" "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
Replace-Text "@@SLOT5@@" "Complement students' training by addressing, in greater depth, current and relevant topics and updating with state-of-the-art topics."
Replace-Text "@@SLOT6@@" "A média do semestre será computada com base na relação: M=(P1+2P2)/3"
Replace-Text "@@SLOT7@@" "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será computada com base na relação: MF=(M+RC)/2"
Replace-Text "@@SLOT8@@" "Apostila ou texto fornecido pelo docente responsável. Artigos extraídos de revistas especializadas nas áreas de Ciências e Tecnologia."
Replace-Text "@@SLOT9@@" "5840730 - Antonio Jefferson da Silva Machado"
Replace-Text "@@SLOT10@@" "519033 - Carlos Yujiro Shigue"

Write-Output "done"
